$wb = $excel.ActiveWorkbook

# Sheet 1: VENTAS POR GRUPO
$ws1 = $wb.Worksheets.Item("VENTAS POR GRUPO")
$ws1.Range("M7").Value = -21.92

# Sheet 2: VENTA MENSUAL
$ws2 = $wb.Worksheets.Item("VENTA MENSUAL")
$ws2.Range("F7").Value = 89.44
$ws2.Range("F22").Value = 2896.86

# Sheet 3: CUMPLIMIENTO MENSUAL
$ws3 = $wb.Worksheets.Item("CUMPLIMIENTO MENSUAL")
$ws3.Range("D16").Value = 2655.9
$ws3.Range("E16").Value = 36120.57
$ws3.Range("F16").Value = 0.06849256778659842

$ws3.Range("D19").Value = 2896.86
$ws3.Range("E19").Value = 56491.36762291769
$ws3.Range("F19").Value = 0.04877835416125659

$ws3.Columns.Item(4).ColumnWidth = 11.166666666666666
